$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 90, shifting rows 90-136 down to 91-137
$ws.Rows.Item(90).Insert()

# Populate the new row 90 with data
$ws.Cells.Item(90, 1).Value = 6
$ws.Cells.Item(90, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(90, 3).Value = "Metropolitana"
$ws.Cells.Item(90, 4).Value = 44529
$ws.Cells.Item(90, 5).Value = 13
$ws.Cells.Item(90, 6).Value = 100112001
$ws.Cells.Item(90, 7).Value = "Berenjena"
$ws.Cells.Item(90, 8).Value = "Sin especificar"
$ws.Cells.Item(90, 9).Value = "Primera"
$ws.Cells.Item(90, 10).Value = 180
$ws.Cells.Item(90, 11).Value = 10000
$ws.Cells.Item(90, 12).Value = 90000
$ws.Cells.Item(90, 13).Value = 45556
$ws.Cells.Item(90, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(90, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(90, 16).Value = 911
$ws.Cells.Item(90, 17).Value = 50
$ws.Cells.Item(90, 18).Value = "Hortaliza"
